$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '52.239.94'
$ws.Range('E2').Value = '  +1.29%  '

# Row 3
$ws.Range('D3').Value = '2.808.45'
$ws.Range('E3').Value = '  +2.22%  '

# Row 4
$ws.Range('E4').Value = '  +0.04%  '

# Row 5
$ws.Range('D5').Value = '347.98'
$ws.Range('E5').Value = '  +4.36%  '

# Row 6
$ws.Range('D6').Value = '115.85'
$ws.Range('E6').Value = '  -0.38%  '

# Row 7
$ws.Range('E7').Value = '  +3.31%  '

# Row 8
$ws.Range('E8').Value = '  -0.03%  '

# Row 9
$ws.Range('D9').Value = '0.591'
$ws.Range('E9').Value = '  +2.33%  '

# Row 10
$ws.Range('D10').Value = '42.38'
$ws.Range('E10').Value = '  +2.30%  '

# Row 11
$ws.Range('D11').Value = '0.0862'
$ws.Range('E11').Value = '  +4.02%  '

# Row 12
$ws.Range('D12').Value = '20.03'
$ws.Range('E12').Value = '  -0.68%  '

# Row 13
$ws.Range('E13').Value = '  +1.44%  '

# Row 14
$ws.Range('D14').Value = '7.84'
$ws.Range('E14').Value = '  +2.92%  '

# Row 15
$ws.Range('D15').Value = '3.245.20'
$ws.Range('E15').Value = '  +2.26%  '

# Row 16
$ws.Range('D16').Value = '2.800.00'
$ws.Range('E16').Value = '  +2.03%  '

# Row 17
$ws.Range('D17').Value = '0.892'
$ws.Range('E17').Value = '  +0.57%  '

# Row 18
$ws.Range('D18').Value = '52.220.96'
$ws.Range('E18').Value = '  +1.40%  '

# Row 19
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').Value = '7.32'
$ws.Range('E19').Value = '  +6.78%  '

# Row 20
$ws.Range('B20').Value = 'ImmutableX'
$ws.Range('C20').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D20').Value = '3.17'
$ws.Range('E20').Value = '  +6.07%  '

# Row 21
$ws.Range('E21').Value = '  -3.04%  '

# Row 22
$ws.Range('D22').Value = '0.0₃0982'
$ws.Range('E22').Value = '  +2.15%  '

# Row 23
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').Value = '70.16'
$ws.Range('E23').Value = '  -0.12%  '

# Row 24
$ws.Range('B24').Value = 'BitcoinCash'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D24').Value = '270.08'
$ws.Range('E24').Value = '  -2.49%  '

# Row 25
$ws.Range('D25').Value = '2.76'
$ws.Range('E25').Value = '  +3.20%  '

# Row 26
$ws.Range('D26').Value = '26.90'
$ws.Range('E26').Value = '  +0.00%  '

# Row 27
$ws.Range('E27').Value = '  -0.06%  '

# Row 28
$ws.Range('D28').Value = '10.28'
$ws.Range('E28').Value = '  -0.19%  '

# Row 29
$ws.Range('E29').Value = '  +1.08%  '

# Row 30
$ws.Range('E30').Value = '  +0.44%  '

# Row 31
$ws.Range('D31').Value = '34.49'
$ws.Range('E31').Value = '  -3.01%  '

# Row 32
$ws.Range('D32').Value = '50.38'
$ws.Range('E32').Value = '  +0.10%  '

# Row 33
$ws.Range('B33').Value = 'VeChain'
$ws.Range('C33').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D33').Value = '0.0445'
$ws.Range('E33').Value = '  +27.68%  '

# Row 34
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '5.78'
$ws.Range('E34').Value = '  +2.91%  '

# Row 35
$ws.Range('D35').Value = '0.0835'
$ws.Range('E35').Value = '  +1.45%  '

# Row 36
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = '2.11'
$ws.Range('E36').Value = '  +0.26%  '

# Row 37
$ws.Range('B37').Value = 'FirstDigitalUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.04%  '

# Row 38
$ws.Range('D38').Value = '18.72'
$ws.Range('E38').Value = '  -3.42%  '

# Row 39
$ws.Range('D39').Value = '4.94'
$ws.Range('E39').Value = '  -1.01%  '

# Row 40
$ws.Range('D40').Value = '3.24'
$ws.Range('E40').Value = '  -1.89%  '

# Row 41
$ws.Range('D41').Value = '2.59'
$ws.Range('E41').Value = '  +9.43%  '

# Row 42
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = '23.42'
$ws.Range('E42').Value = '  -0.92%  '

# Row 43
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').Value = '0.116'
$ws.Range('E43').Value = '  +1.98%  '

# Row 44
$ws.Range('D44').Value = '126.40'
$ws.Range('E44').Value = '  -2.40%  '

# Row 45
$ws.Range('E45').Value = '  +0.12%  '

# Row 46
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').Value = '3.32'
$ws.Range('E46').Value = '  -1.71%  '

# Row 47
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '2.060.64'
$ws.Range('E47').Value = '  -1.93%  '

# Row 48
$ws.Range('E48').Value = '  +3.31%  '

# Row 49
$ws.Range('D49').Value = '0.967'
$ws.Range('E49').Value = '  +12.99%  '

# Row 50
$ws.Range('D50').Value = '5.60'
$ws.Range('E50').Value = '  -0.44%  '

# Row 51
$ws.Range('D51').Value = '9.00'
$ws.Range('E51').Value = '  +0.33%  '
